$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated loading_percent values (case with 380 kV)
$ws.Range("B2").Value = 16.72772944972617
$ws.Range("C2").Value = 9.985852220438849
$ws.Range("D2").Value = 10.08902423161724
$ws.Range("F2").Value = 35.24802390981363
$ws.Range("G2").Value = 37.93853244728479
$ws.Range("H2").Value = 16.13047118751786
$ws.Range("J2").Value = 11.35377249047149
$ws.Range("N2").Value = 17.21177631290732

$ws.Range("B3").Value = 16.16419213718873
$ws.Range("C3").Value = 9.456785351342521
$ws.Range("D3").Value = 10.04248633092132
$ws.Range("F3").Value = 35.11769566837565
$ws.Range("G3").Value = 37.64632634482904
$ws.Range("H3").Value = 16.15305342033729
$ws.Range("J3").Value = 11.33834666643095
$ws.Range("N3").Value = 17.28317986346202

$ws.Range("B4").Value = 15.81177377851516
$ws.Range("C4").Value = 9.117763237468914
$ws.Range("D4").Value = 10.01549962023959
$ws.Range("F4").Value = 35.04914211878977
$ws.Range("G4").Value = 37.48171123223569
$ws.Range("H4").Value = 16.17106108166537
$ws.Range("J4").Value = 11.33131085130395
$ws.Range("N4").Value = 17.32899528429521

$ws.Range("B5").Value = 15.66678817205626
$ws.Range("C5").Value = 8.976200238910327
$ws.Range("D5").Value = 10.00490913525305
$ws.Range("F5").Value = 35.02410464468846
$ws.Range("G5").Value = 37.41841188086753
$ws.Range("H5").Value = 16.17943708146122
$ws.Range("J5").Value = 11.32905749728141
$ws.Range("N5").Value = 17.34816345564698

$ws.Range("B6").Value = 15.6426378752471
$ws.Range("C6").Value = 8.952493025695674
$ws.Range("D6").Value = 10.00317537540994
$ws.Range("F6").Value = 35.02012261030301
$ws.Range("G6").Value = 37.40813113820916
$ws.Range("H6").Value = 16.18089046665062
$ws.Range("J6").Value = 11.32872043100537
$ws.Range("N6").Value = 17.35137644774478

$ws.Range("B7").Value = 15.80982367947306
$ws.Range("C7").Value = 9.115867645381002
$ws.Range("D7").Value = 10.01535513644126
$ws.Range("F7").Value = 35.04879270193796
$ws.Range("G7").Value = 37.48084216561129
$ws.Range("H7").Value = 16.1711698470251
$ws.Range("J7").Value = 11.33127797514169
$ws.Range("N7").Value = 17.32925177416443

$ws.Range("B8").Value = 16.53488784927742
$ws.Range("C8").Value = 9.806456224228429
$ws.Range("D8").Value = 10.07265311220637
$ws.Range("F8").Value = 35.20071662432778
$ws.Range("G8").Value = 37.83474717030143
$ws.Range("H8").Value = 16.13739566579194
$ws.Range("J8").Value = 11.34794906143689
$ws.Range("N8").Value = 17.23598777481076

$ws.Range("B9").Value = 17.89645004770693
$ws.Range("C9").Value = 11.04273219416118
$ws.Range("D9").Value = 10.19725872496537
$ws.Range("F9").Value = 35.5887857608389
$ws.Range("G9").Value = 38.64308778720473
$ws.Range("H9").Value = 16.10418506618468
$ws.Range("J9").Value = 11.3998906377687
$ws.Range("N9").Value = 17.06867196674729

$ws.Range("B10").Value = 18.8487391983872
$ws.Range("C10").Value = 11.87325818788558
$ws.Range("D10").Value = 10.29576862791191
$ws.Range("F10").Value = 35.92743453573257
$ws.Range("G10").Value = 39.3019650641531
$ws.Range("H10").Value = 16.10010290313529
$ws.Range("J10").Value = 11.44965759956261
$ws.Range("N10").Value = 16.95512216609556

$ws.Range("B11").Value = 19.26952278614161
$ws.Range("C11").Value = 12.23332279584763
$ws.Range("D11").Value = 10.3419796037071
$ws.Range("F11").Value = 36.09273174936655
$ws.Range("G11").Value = 39.61467164646722
$ws.Range("H11").Value = 16.10268855123218
$ws.Range("J11").Value = 11.47478071854703
$ws.Range("N11").Value = 16.90547649872824

$ws.Range("B12").Value = 19.42693435493752
$ws.Range("C12").Value = 12.36706322964205
$ws.Range("D12").Value = 10.35966925569578
$ws.Range("F12").Value = 36.15690394365798
$ws.Range("G12").Value = 39.73484791308572
$ws.Range("H12").Value = 16.10430821661344
$ws.Range("J12").Value = 11.48464743077296
$ws.Range("N12").Value = 16.88696395260752

$ws.Range("B13").Value = 19.39312091782591
$ws.Range("C13").Value = 12.33837661306669
$ws.Range("D13").Value = 10.35585116818146
$ws.Range("F13").Value = 36.1430138040066
$ws.Range("G13").Value = 39.70888920291265
$ws.Range("H13").Value = 16.10393088305638
$ws.Range("J13").Value = 11.4825068213302
$ws.Range("N13").Value = 16.89093821401362

$ws.Range("B14").Value = 19.28251257914646
$ws.Range("C14").Value = 12.24437825695039
$ws.Range("D14").Value = 10.34343118427272
$ws.Range("F14").Value = 36.09797985109991
$ws.Range("G14").Value = 39.6245239446616
$ws.Range("H14").Value = 16.1028089568746
$ws.Range("J14").Value = 11.47558540638374
$ws.Range("N14").Value = 16.90394771414162

$ws.Range("B15").Value = 19.21450643043869
$ws.Range("C15").Value = 12.18646040458632
$ws.Range("D15").Value = 10.33584808445969
$ws.Range("F15").Value = 36.07059949933916
$ws.Range("G15").Value = 39.57307392957541
$ws.Range("H15").Value = 16.10220520307075
$ws.Range("J15").Value = 11.47139170032318
$ws.Range("N15").Value = 16.91195376032169

$ws.Range("B16").Value = 18.82097739604251
$ws.Range("C16").Value = 11.84936494002088
$ws.Range("D16").Value = 10.29277587772835
$ws.Range("F16").Value = 35.91685496738973
$ws.Range("G16").Value = 39.2817812314221
$ws.Range("H16").Value = 16.10002347812964
$ws.Range("J16").Value = 11.44806539111064
$ws.Range("N16").Value = 16.95840690537953

$ws.Range("B17").Value = 18.57627147109734
$ws.Range("C17").Value = 11.63797875625214
$ws.Range("D17").Value = 10.26670323455041
$ws.Range("F17").Value = 35.82539083580478
$ws.Range("G17").Value = 39.1063334955892
$ws.Range("H17").Value = 16.09982428831157
$ws.Range("J17").Value = 11.434388902449
$ws.Range("N17").Value = 16.98741769811426

$ws.Range("B18").Value = 18.43436181254638
$ws.Range("C18").Value = 11.51472724552051
$ws.Range("D18").Value = 10.25183920755761
$ws.Range("F18").Value = 35.7738440199344
$ws.Range("G18").Value = 39.00664896214514
$ws.Range("H18").Value = 16.10012781271278
$ws.Range("J18").Value = 11.42675661377045
$ws.Range("N18").Value = 17.00429310605618

$ws.Range("B19").Value = 18.38611886871089
$ws.Range("C19").Value = 11.47271187231099
$ws.Range("D19").Value = 10.24682953592223
$ws.Range("F19").Value = 35.75657450577244
$ws.Range("G19").Value = 38.97311160425471
$ws.Range("H19").Value = 16.10030232905702
$ws.Range("J19").Value = 11.42421276882382
$ws.Range("N19").Value = 17.01003936963707

$ws.Range("B20").Value = 18.60244212283764
$ws.Range("C20").Value = 11.66065422925726
$ws.Range("D20").Value = 10.26946509794891
$ws.Range("F20").Value = 35.83501781912784
$ws.Range("G20").Value = 39.12488382591279
$ws.Range("H20").Value = 16.09980220770287
$ws.Range("J20").Value = 11.43582059097966
$ws.Range("N20").Value = 16.98430988376951

$ws.Range("B21").Value = 19.3150543779577
$ws.Range("C21").Value = 12.27205903196629
$ws.Range("D21").Value = 10.34707414823951
$ws.Range("F21").Value = 36.11116492498919
$ws.Range("G21").Value = 39.64925713223308
$ws.Range("H21").Value = 16.10312109882019
$ws.Range("J21").Value = 11.47760884430733
$ws.Range("N21").Value = 16.90011872747272

$ws.Range("B22").Value = 19.7694825158696
$ws.Range("C22").Value = 12.65642621845792
$ws.Range("D22").Value = 10.39890178633534
$ws.Range("F22").Value = 36.3008186699622
$ws.Range("G22").Value = 40.00217860188818
$ws.Range("H22").Value = 16.10902415235341
$ws.Range("J22").Value = 11.50697601640238
$ws.Range("N22").Value = 16.84676814738386

$ws.Range("B23").Value = 19.52802351820824
$ws.Range("C23").Value = 12.45269072651052
$ws.Range("D23").Value = 10.37114278551223
$ws.Range("F23").Value = 36.19877114625167
$ws.Range("G23").Value = 39.81291913570344
$ws.Range("H23").Value = 16.10553149471943
$ws.Range("J23").Value = 11.49111552035512
$ws.Range("N23").Value = 16.87508980171643

$ws.Range("B24").Value = 18.5906141755069
$ws.Range("C24").Value = 11.65040800963753
$ws.Range("D24").Value = 10.26821606785559
$ws.Range("F24").Value = 35.83066222644545
$ws.Range("G24").Value = 39.11649352958842
$ws.Range("H24").Value = 16.09981088828935
$ws.Range("J24").Value = 11.43517260630597
$ws.Range("N24").Value = 16.98571431373066

$ws.Range("B25").Value = 17.53580762720258
$ws.Range("C25").Value = 10.72159429162739
$ws.Range("D25").Value = 10.16228890240512
$ws.Range("F25").Value = 35.47428236751271
$ws.Range("G25").Value = 38.41262797761942
$ws.Range("H25").Value = 16.10961321782586
$ws.Range("J25").Value = 11.38378977136335
$ws.Range("N25").Value = 17.11228016542953
